$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: add B19 "10 uur 15 minuten" (bold, like A19/B17) ---
$ws.Range("B19").Value = "10 uur 15 minuten"
$ws.Range("B19").Font.Bold = $true

# --- Row 22: add A22 "Week 7" (bold, like A19) ---
$ws.Range("A22").Value = "Week 7"
$ws.Range("A22").Font.Bold = $true

# --- Row 23: new log entry ---
$ws.Range("A23").Value = 43553
$ws.Range("A23").NumberFormat = $ws.Range("A21").NumberFormat

$ws.Range("B23").Value = "4 uur"

$ws.Range("C23").Value = "maps fragment aanmaken met functionaliteit routes maken door markers toe te voegen, fragment info activiteit toevoegen met mapview en ercycler view, alle navigatie toevoegen en codelab friendlychat project maken"

# --- Row 24: new hyperlink (stackoverflow) ---
$ws.Hyperlinks.Add($ws.Range("Q24"), "https://stackoverflow.com/questions/13932441/android-google-maps-v2-set-zoom-level-for-mylocation")
$ws.Range("Q24").Style = $ws.Range("Q22").Style

# --- Row 25: new hyperlink (codelabs, with location/subaddress "0") ---
$ws.Hyperlinks.Add($ws.Range("Q25"), "https://codelabs.developers.google.com/codelabs/firebase-android/", "0", "", "https://codelabs.developers.google.com/codelabs/firebase-android/ - 0")
$ws.Range("Q25").Value = "https://codelabs.developers.google.com/codelabs/firebase-android/#0"
$ws.Range("Q25").Style = $ws.Range("Q22").Style

# --- Row 26: new hyperlink (firebase console) ---
$ws.Hyperlinks.Add($ws.Range("Q26"), "https://console.firebase.google.com/project/friendlychat-24131/database/friendlychat-24131/rules")
$ws.Range("Q26").Style = $ws.Range("Q22").Style

# --- Selection change ---
[void]$ws.Range("G30").Select()
